$wb = $excel.ActiveWorkbook

# --- TTD sheet: add a new "Price Type" column and a new test segment row ---
$ws = $wb.Worksheets.Item("TTD")

# Insert a new column before the existing "Batch ID" column (I), shifting
# Batch ID / Depth / Segment Full Path one column to the right.
[void]$ws.Columns("I").Insert()
$ws.Columns("I").ColumnWidth = $ws.Columns("H").ColumnWidth

# Header + helper-text rows for the new "Price Type" column.
$ws.Range("I1").Value = "Price Type"
$ws.Range("I2").Value = "Add/Edit: Required`nEdit Rates: Required`nRetrieve Batch: Optional`nRetrieve Rates: Optional`nValues: CPM or PercentOfMediaCost"
$ws.Rows("2:2").RowHeight = 102
$ws.Range("I3").Value = "CPM"

# New data row (order matters so new shared strings are appended in the
# same sequence as the authored workbook).
$ws.Range("A4").Value = 20190401004
$ws.Range("B4").Value = "ttdratetest_partnerID_rate"
$ws.Range("I4").Value = "PercentOfMediaCost"
$ws.Range("F4").Value = "eyeota"
$ws.Range("C4").Value = "Test Segment 20190401004"
$ws.Range("D4").Value = "Test Segment 20190401004"
$ws.Range("E4").Value = $true
$ws.Range("G4").Value = "abc123"
$ws.Range("H4").Value = 1
$ws.Range("J4").Value = 23456
$ws.Range("H4:J4").Style = "Normal"

# --- Selection / active-tab bookkeeping ---
$wsAdform = $wb.Worksheets.Item("Adform")
[void]$wsAdform.Activate()
[void]$wsAdform.Range("C24").Select()

[void]$ws.Activate()
[void]$ws.Range("B5").Select()
